$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '30.287.01'
$dCell.Style = "Normal"
$eCell = $ws.Range("E2")
$eCell.NumberFormat = "@"
$eCell.Value = '  +5.31%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '1.921.61'
$dCell.Style = "Normal"
$eCell = $ws.Range("E3")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.28%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9968'
$dCell.Style = "Normal"
$eCell = $ws.Range("E4")
$eCell.NumberFormat = "@"
$eCell.Value = '  -0.30%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '253.58'
$dCell.Style = "Normal"
$eCell = $ws.Range("E5")
$eCell.NumberFormat = "@"
$eCell.Value = '  +0.29%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9972'
$dCell.Style = "Normal"
$eCell = $ws.Range("E6")
$eCell.NumberFormat = "@"
$eCell.Value = '  -0.23%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = '0.5195'
$dCell.Style = "Normal"
$eCell = $ws.Range("E7")
$eCell.NumberFormat = "@"
$eCell.Value = '  +4.69%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '46.12'
$dCell.Style = "Normal"
$eCell = $ws.Range("E8")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.87%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = '0.2992'
$dCell.Style = "Normal"
$eCell = $ws.Range("E9")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.17%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = '0.06813'
$dCell.Style = "Normal"
$eCell = $ws.Range("E10")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.68%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '1.915.20'
$dCell.Style = "Normal"
$eCell = $ws.Range("E11")
$eCell.NumberFormat = "@"
$eCell.Value = '  +5.91%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = '17.61'
$dCell.Style = "Normal"
$eCell = $ws.Range("E12")
$eCell.NumberFormat = "@"
$eCell.Value = '  +4.27%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '0.07323'
$dCell.Style = "Normal"
$eCell = $ws.Range("E13")
$eCell.NumberFormat = "@"
$eCell.Value = '  +3.00%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = '0.6908'
$dCell.Style = "Normal"
$eCell = $ws.Range("E14")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.84%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '88.30'
$dCell.Style = "Normal"
$eCell = $ws.Range("E15")
$eCell.NumberFormat = "@"
$eCell.Value = '  +7.76%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '4.938'
$dCell.Style = "Normal"
$eCell = $ws.Range("E16")
$eCell.NumberFormat = "@"
$eCell.Value = '  +4.93%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '30.269.04'
$dCell.Style = "Normal"
$eCell = $ws.Range("E17")
$eCell.NumberFormat = "@"
$eCell.Value = '  +5.31%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = '0.000007828'
$dCell.Style = "Normal"
$eCell = $ws.Range("E18")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.44%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9975'
$dCell.Style = "Normal"
$eCell = $ws.Range("E19")
$eCell.NumberFormat = "@"
$eCell.Value = '  -0.19%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = '13.14'
$dCell.Style = "Normal"
$eCell = $ws.Range("E20")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.96%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '2.160.14'
$dCell.Style = "Normal"
$eCell = $ws.Range("E21")
$eCell.NumberFormat = "@"
$eCell.Value = '  +5.86%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9965'
$dCell.Style = "Normal"
$eCell = $ws.Range("E22")
$eCell.NumberFormat = "@"
$eCell.Value = '  -0.31%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '4.874'
$dCell.Style = "Normal"
$eCell = $ws.Range("E23")
$eCell.NumberFormat = "@"
$eCell.Value = '  +5.65%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '5.761'
$dCell.Style = "Normal"
$eCell = $ws.Range("E24")
$eCell.NumberFormat = "@"
$eCell.Value = '  +8.73%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '9.252'
$dCell.Style = "Normal"
$eCell = $ws.Range("E25")
$eCell.NumberFormat = "@"
$eCell.Value = '  +4.00%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = '139.80'
$dCell.Style = "Normal"
$eCell = $ws.Range("E26")
$eCell.NumberFormat = "@"
$eCell.Value = '  +25.27%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = '146.41'
$dCell.Style = "Normal"
$eCell = $ws.Range("E27")
$eCell.NumberFormat = "@"
$eCell.Value = '  +2.70%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '17.43'
$dCell.Style = "Normal"
$eCell = $ws.Range("E28")
$eCell.NumberFormat = "@"
$eCell.Value = '  +8.84%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = '2.025'
$dCell.Style = "Normal"
$eCell = $ws.Range("E29")
$eCell.NumberFormat = "@"
$eCell.Value = '  +7.49%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = '1.382'
$dCell.Style = "Normal"
$eCell = $ws.Range("E30")
$eCell.NumberFormat = "@"
$eCell.Value = '  -0.53%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '4.301'
$dCell.Style = "Normal"
$eCell = $ws.Range("E31")
$eCell.NumberFormat = "@"
$eCell.Value = '  +2.87%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = '0.08895'
$dCell.Style = "Normal"
$eCell = $ws.Range("E32")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.31%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = '4.061'
$dCell.Style = "Normal"
$eCell = $ws.Range("E33")
$eCell.NumberFormat = "@"
$eCell.Value = '  +5.45%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '0.05161'
$dCell.Style = "Normal"
$eCell = $ws.Range("E34")
$eCell.NumberFormat = "@"
$eCell.Value = '  +3.96%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = '1.164'
$dCell.Style = "Normal"
$eCell = $ws.Range("E35")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.25%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '0.7221'
$dCell.Style = "Normal"
$eCell = $ws.Range("E36")
$eCell.NumberFormat = "@"
$eCell.Value = '  +7.42%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = '2.687'
$dCell.Style = "Normal"
$eCell = $ws.Range("E37")
$eCell.NumberFormat = "@"
$eCell.Value = '  +0.65%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '2.862'
$dCell.Style = "Normal"
$eCell = $ws.Range("E38")
$eCell.NumberFormat = "@"
$eCell.Value = '  +8.64%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = '2.327'
$dCell.Style = "Normal"
$eCell = $ws.Range("E39")
$eCell.NumberFormat = "@"
$eCell.Value = '  +8.01%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9731'
$dCell.Style = "Normal"
$eCell = $ws.Range("E40")
$eCell.NumberFormat = "@"
$eCell.Value = '  +0.98%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '0.01699'
$dCell.Style = "Normal"
$eCell = $ws.Range("E41")
$eCell.NumberFormat = "@"
$eCell.Value = '  +5.75%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '6.161'
$dCell.Style = "Normal"
$eCell = $ws.Range("E42")
$eCell.NumberFormat = "@"
$eCell.Value = '  +3.45%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '0.4364'
$dCell.Style = "Normal"
$eCell = $ws.Range("E43")
$eCell.NumberFormat = "@"
$eCell.Value = '  +5.75%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '106.26'
$dCell.Style = "Normal"
$eCell = $ws.Range("E44")
$eCell.NumberFormat = "@"
$eCell.Value = '  +5.30%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9989'
$dCell.Style = "Normal"
$eCell = $ws.Range("E45")
$eCell.NumberFormat = "@"
$eCell.Value = '  -0.01%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = '7.702'
$dCell.Style = "Normal"
$eCell = $ws.Range("E46")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.48%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = '0.1283'
$dCell.Style = "Normal"
$eCell = $ws.Range("E47")
$eCell.NumberFormat = "@"
$eCell.Value = '  +4.51%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = '0.05732'
$dCell.Style = "Normal"
$eCell = $ws.Range("E48")
$eCell.NumberFormat = "@"
$eCell.Value = '  +4.31%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = '8.565'
$dCell.Style = "Normal"
$eCell = $ws.Range("E49")
$eCell.NumberFormat = "@"
$eCell.Value = '  +4.24%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '33.26'
$dCell.Style = "Normal"
$eCell = $ws.Range("E50")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.22%  '
$eCell.Style = "Normal"
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = '0.3868'
$dCell.Style = "Normal"
$eCell = $ws.Range("E51")
$eCell.NumberFormat = "@"
$eCell.Value = '  +6.81%  '
$eCell.Style = "Normal"
